# Update Nova autumn (2024 Q4) character nominations
# - Fix stray number-format style on the final pre-existing row (F768:G768)
# - Append 31 new character rows (769-799) pulled from the autumn nomination list
# - Refresh the sheet selection to match the author's final cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix number-format style drift on row 768 (was s="5", should be s="1") ---
$ws.Range("F765:G765").Copy()
$ws.Range("F768:G768").PasteSpecial(-4122)

# --- 2. Append the new rows starting at row 769 ---
$rows = @(
    @("猪股大喜", "青之箱", "千叶翔也", "Inomata Taiki", 2024, 10),
    @("高仓健", "胆大党", "花江夏树", "Takakura Ken", 2024, 10),
    @("田沼要", "夏目友人帐", "堀江一真", "Tanuma Kaname", 2008, 7),
    @("鸭乃桥论", "鸭乃桥论的禁忌推理", "阿座上洋平", "Kamonohashi Ron", 2023, 10),
    @("藤丸立香", "Fate系列", "岛崎信长", "Fujimaru Ritsuka", 2006, 1),
    @("早乙女乱马", "乱马½", "山口胜平", "Saotome Ranma", 1989, 4),
    @("中野丸尾", "五等分的新娘", "黑田崇矢", "Nakano Maruo", 2019, 1),
    @("多兰", "再见龙生，你好人生", "武内骏辅", "Dolan", 2024, 10),
    @("佐佐木常宏", "悲喜渔生", "岩中睦树", "Sasaki Tsunehiro", 2024, 10),
    @("拉法尔", "地。-关于地球的运动-", "坂本真绫", "Rafal", 2024, 10),
    @("圆城寺仁", "胆大党", "石川界人", "Enjōji Jin", 2024, 10),
    @("大原拓也", "听说你们要结婚！？", "熊谷健太郎", "Ōhara Takuya", 2024, 10),
    @("贝尔多尔·贝尔别特·贝尔休伯特", "魔王2099", "日野聪", "Veltol Velvet Velsvalt", 2024, 10),
    @("上终瓜生", "缘结甘神家", "铃木崚汰", "Kamihate Uryū", 2024, 10),
    @("踯躅森贵明", "悲喜渔生", "石川界人", "Tsutsujimori Takaaki", 2024, 10),
    @("花散仁央", "青之壬生浪", "梅田修一朗", "Chirinu Nio", 2024, 10),
    @("哈迪斯·迪奥斯·拉维", "重启人生的千金小姐正在攻略龙帝陛下", "户谷菊之介", "Hadis Teos Rave", 2024, 10),
    @("河合井小太郎", "噗妮露是可爱史莱姆", "梅田修一朗", "Kawaii Kotarō", 2024, 10),
    @("克莱·安东黎希", "叹气的亡灵想隐退", "小野贤章", "Krai Andrey", 2024, 10),
    @("新岛圭介", "妻子变成小学生。", "平川大辅", "Niijima Keisuke", 2024, 10),
    @("南云始", "平凡职业成就世界最强", "深町寿成", "Nagumo Hajime", 2019, 7),
    @("天束光", "机械手臂", "丰永利行", "Amatsuga Hikaru", 2024, 10),
    @("库洛马", "Acro Trip 顶尖恶路", "岛崎信长", "Chrome", 2024, 10),
    @("遥", "孤单一人的异世界攻略", "梅田修一朗", "Haruka", 2024, 10),
    @("阳务乐郎", "香格里拉边境～粪作猎人向神作游戏发起挑战～", "内田雄马", "Hizutome Rakurō", 2023, 10),
    @("五十岚一贺", "喂！蜻蜓", "东地宏树", "Igarashi Kazuyoshi", 2024, 4),
    @("雷格西", "BEASTARS", "小林亲弘", "Legoshi", 2019, 10),
    @("奥托·苏文", "Re:从零开始的异世界生活", "天崎滉平", "Otto Suwen", 2016, 4),
    @("利欧", "精灵幻想记", "松冈祯丞", "Rio", 2021, 7),
    @("重本浩司", "魔法光源股份有限公司", "小山力也", "Shigemoto Kōji", 2024, 10),
    @("常盘", "去参加联谊，却发现完全没有女生在场", "武内骏辅", "Tokiwa", 2024, 10)
)

$startRow = 769
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
    $ws.Range("D$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}

$endRow = $startRow + $rows.Count - 1

# Row height to match the rest of the table (20pt, custom)
$ws.Rows("$startRow`:$endRow").RowHeight = 20

# Copy the plain number format (style s="1") onto the new F:G cells
$ws.Range("F765:G765").Copy()
$ws.Range("F$startRow`:G$endRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3. Restore the author's final selection ---
$ws.Range("E793").Select()
